# Changes of 27th April 2022
#
# The underlying MDSi test-result sheet records the job numbers that were
# processed in the most recent test run. B2:B4 need to be updated to the
# job numbers from the latest run while keeping them as text values
# (shared strings), matching their original "t=s" cell type and their
# original (unstyled) cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force a Text number format so the numeric-looking job
# numbers are stored as text (shared strings) rather than being
# auto-coerced to numeric cells -- this mirrors how the original workbook
# stores B2:B4 (t="s").
$ws.Range("B2:B4").NumberFormat = "@"

$ws.Range("B2").Value = "32372249"
$ws.Range("B3").Value = "32372251"
$ws.Range("B4").Value = "32372252"

# Restore the default ("Normal") style so the cells keep their original
# (unstyled) appearance, same as before the edit.
$ws.Range("B2:B4").Style = "Normal"
